# Add 2022-Q3 data:
#  - Insert a new worksheet "2022-Q3" right after the "总计" summary sheet,
#    carrying the usual per-quarter fund-holdings table for a single fund.
#  - Insert a matching summary row at the top of the "总计" sheet's data
#    (row 2), pushing the existing quarters down and renumbering the index
#    column.

$wb = $excel.ActiveWorkbook

$total = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1) New "2022-Q3" worksheet, inserted right after "总计"
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Add($null, $total)
$q3.Name = "2022-Q3"

$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

# Match the bold / thin-bordered / centered header look used by every
# other per-quarter sheet.
$hdr = $q3.Range("B1:H1")
$hdr.Font.Bold = $true
$hdr.Borders.LineStyle = 1
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160

$q3.Range("A2").Value = 0
$q3.Range("A2").Font.Bold = $true
$q3.Range("A2").Borders.LineStyle = 1
$q3.Range("A2").HorizontalAlignment = -4108
$q3.Range("A2").VerticalAlignment = -4160

$q3.Range("B2").NumberFormat = "@"
$q3.Range("B2").Value = "161838"

$q3.Range("C2").Value = "银华创业板两年定期开放混合"

$q3.Range("D2").NumberFormat = "@"
$q3.Range("D2").Value = "4.45"

$q3.Range("E2").NumberFormat = "@"
$q3.Range("E2").Value = "96.33"

$q3.Range("F2").NumberFormat = "@"
$q3.Range("F2").Value = "8.84"

$q3.Range("G2").NumberFormat = "@"
$q3.Range("G2").Value = "0.3934"

$q3.Range("H2").Value = 6

# ---------------------------------------------------------------------
# 2) "总计" sheet: insert a new row 2 for 2022-Q3, push the rest down,
#    and renumber the index column (A) to stay 0-based sequential.
# ---------------------------------------------------------------------
$total.Rows.Item(2).Insert()

# The inserted blank row doesn't automatically pick up the bordered /
# bold index-column style used by the rest of column A - copy it down
# from the row just below (which still has the original formatting).
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 1
$total.Range("D2").Value = 0.39

$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
$total.Range("A7").Value = 5

# Restore the originally-selected tab ("2020-Q4" was the active sheet
# before this edit) so the insert/activate calls above don't leave an
# unrelated side effect on the saved selection state.
$wb.Worksheets.Item("2020-Q4").Activate()

